# Update crypto price/volume snapshot values (D: Price, E: Volume(1h)).
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as text (matching the source data, e.g. "5.24", "1.00"),
# then the style is reset to "Normal" so no explicit cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.309.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "'2.433.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'563.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "

$ws.Range("D6").Value = "'144.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").Value = "'2.429.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "'5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").Value = "'0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("D14").Value = "'26.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.08%  "

$ws.Range("D15").Value = "'0.0000175"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.37%  "

$ws.Range("D16").Value = "'2.859.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").Value = "'62.269.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").Value = "'2.434.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "'11.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").Value = "'325.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("E21").Value = "  -1.32%  "

$ws.Range("D22").Value = "'6.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'67.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.52%  "

$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "'8.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "

$ws.Range("D27").Value = "'554.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.42%  "

$ws.Range("D28").Value = "'2.544.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "'0.0₃0943"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("D31").Value = "'8.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "'1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.52%  "

$ws.Range("D33").Value = "'0.148"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").Value = "'4.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "

$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").Value = "'5.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "

$ws.Range("D40").Value = "'18.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").Value = "'150.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.14%  "

$ws.Range("D42").Value = "'1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.58%  "

$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").Value = "'2.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.16%  "

$ws.Range("D45").Value = "'148.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "

$ws.Range("D46").Value = "'3.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("D47").Value = "'0.0535"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Value = "'20.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("E51").Value = "  +0.40%  "
